$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Cells.Item(3, 8).Value = 19000  # H3
$ws.Cells.Item(3, 10).Value = 19000  # J3
$ws.Cells.Item(3, 12).Value = 19000  # L3
$ws.Cells.Item(3, 14).Value = -19228  # N3
# Row 38
$ws.Cells.Item(38, 8).Value = 485.5  # H38
$ws.Cells.Item(38, 9).Value = 197.85715  # I38
$ws.Cells.Item(38, 11).Value = 593.5714499999999  # K38
$ws.Cells.Item(38, 13).Value = -221.5714499999999  # M38
# Row 51
$ws.Cells.Item(51, 8).Value = 50103900  # H51
$ws.Cells.Item(51, 9).Value = 171666.33  # I51
$ws.Cells.Item(51, 11).Value = 171666.33  # K51
$ws.Cells.Item(51, 13).Value = -171182.33  # M51
# Row 76
$ws.Cells.Item(76, 8).Value = 4949.5  # H76
$ws.Cells.Item(76, 9).Value = 4949.5  # I76
$ws.Cells.Item(76, 11).Value = 4949.5  # K76
$ws.Cells.Item(76, 13).Value = -4634.5  # M76
# Row 79
$ws.Cells.Item(79, 8).Value = 4949.5  # H79
$ws.Cells.Item(79, 9).Value = 4949.5  # I79
$ws.Cells.Item(79, 11).Value = 4949.5  # K79
$ws.Cells.Item(79, 13).Value = -3857.5  # M79
# Row 102
$ws.Cells.Item(102, 8).Value = 19000  # H102
$ws.Cells.Item(102, 10).Value = 19000  # J102
$ws.Cells.Item(102, 12).Value = 19000  # L102
$ws.Cells.Item(102, 14).Value = -25490  # N102
# Row 106
$ws.Cells.Item(106, 8).Value = 1589.35  # H106
$ws.Cells.Item(106, 9).Value = 1443.2812  # I106
$ws.Cells.Item(106, 10).Value = 2173.625  # J106
$ws.Cells.Item(106, 11).Value = 1443.2812  # K106
$ws.Cells.Item(106, 12).Value = 2173.625  # L106
$ws.Cells.Item(106, 13).Value = -812.2811999999999  # M106
$ws.Cells.Item(106, 14).Value = -3435.625  # N106
# Row 111
$ws.Cells.Item(111, 8).Value = 7500  # H111
$ws.Cells.Item(111, 10).Value = 0  # J111
$ws.Cells.Item(111, 12).Value = 0  # L111
$ws.Cells.Item(111, 14).ClearContents()  # N111
# Row 112
$ws.Cells.Item(112, 8).Value = 32011.514  # H112
$ws.Cells.Item(112, 9).Value = 2491.2856  # I112
$ws.Cells.Item(112, 10).Value = 39391.57  # J112
$ws.Cells.Item(112, 11).Value = 7473.8568  # K112
$ws.Cells.Item(112, 12).Value = 118174.71  # L112
$ws.Cells.Item(112, 13).Value = -6365.8568  # M112
$ws.Cells.Item(112, 14).Value = -120390.71  # N112
# Row 115
$ws.Cells.Item(115, 8).Value = 821.25  # H115
$ws.Cells.Item(115, 9).Value = 821.25  # I115
$ws.Cells.Item(115, 11).Value = 2463.75  # K115
$ws.Cells.Item(115, 13).Value = -896.75  # M115
# Row 138
$ws.Cells.Item(138, 8).Value = 5145.6724  # H138
$ws.Cells.Item(138, 10).Value = 3582.4814  # J138
$ws.Cells.Item(138, 12).Value = 10747.4442  # L138
$ws.Cells.Item(138, 14).Value = -21027.4442  # N138

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 158399.23  # H32
$ws.Cells.Item(32, 9).Value = 170466.84  # I32
$ws.Cells.Item(32, 11).Value = 170466.84  # K32
$ws.Cells.Item(32, 13).Value = -170179.84  # M32
# Row 63
$ws.Cells.Item(63, 8).Value = 20516.137  # H63
$ws.Cells.Item(63, 9).Value = 5498.3335  # I63
$ws.Cells.Item(63, 11).Value = 5498.3335  # K63
$ws.Cells.Item(63, 13).Value = -4812.3335  # M63
# Row 66
$ws.Cells.Item(66, 8).Value = 20516.137  # H66
$ws.Cells.Item(66, 9).Value = 5498.3335  # I66
$ws.Cells.Item(66, 11).Value = 27491.6675  # K66
$ws.Cells.Item(66, 13).Value = -24059.6675  # M66
# Row 74
$ws.Cells.Item(74, 8).Value = 373004.84  # H74
$ws.Cells.Item(74, 9).Value = 1056.6086  # I74
$ws.Cells.Item(74, 11).Value = 1056.6086  # K74
$ws.Cells.Item(74, 13).Value = -182.6086  # M74
# Row 77
$ws.Cells.Item(77, 8).Value = 373004.84  # H77
$ws.Cells.Item(77, 9).Value = 1056.6086  # I77
$ws.Cells.Item(77, 11).Value = 5283.043  # K77
$ws.Cells.Item(77, 13).Value = -915.0429999999997  # M77
# Row 92
$ws.Cells.Item(92, 8).Value = 47000  # H92
$ws.Cells.Item(92, 10).Value = 47000  # J92
$ws.Cells.Item(92, 12).Value = 47000  # L92
$ws.Cells.Item(92, 14).Value = -51992  # N92
# Row 96
$ws.Cells.Item(96, 8).Value = 49185  # H96
$ws.Cells.Item(96, 10).Value = 49185  # J96
$ws.Cells.Item(96, 12).Value = 49185  # L96
$ws.Cells.Item(96, 14).Value = -54677  # N96
# Row 138
$ws.Cells.Item(138, 8).Value = 95000  # H138
$ws.Cells.Item(138, 10).Value = 95000  # J138
$ws.Cells.Item(138, 12).Value = 95000  # L138
$ws.Cells.Item(138, 14).Value = -105280  # N138

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 9
$ws.Cells.Item(9, 8).Value = 100000  # H9
$ws.Cells.Item(9, 10).Value = 100000  # J9
$ws.Cells.Item(9, 12).Value = 100000  # L9
$ws.Cells.Item(9, 14).Value = -100336  # N9
# Row 20
$ws.Cells.Item(20, 8).Value = 1103.7894  # H20
$ws.Cells.Item(20, 10).Value = 935.8889  # J20
$ws.Cells.Item(20, 12).Value = 935.8889  # L20
$ws.Cells.Item(20, 14).Value = -1429.8889  # N20

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Cells.Item(22, 8).Value = 1003.3214  # H22
$ws.Cells.Item(22, 9).Value = 1007.24  # I22
$ws.Cells.Item(22, 10).Value = 970.6667  # J22
$ws.Cells.Item(22, 11).Value = 1007.24  # K22
$ws.Cells.Item(22, 12).Value = 970.6667  # L22
$ws.Cells.Item(22, 13).Value = -657.24  # M22
$ws.Cells.Item(22, 14).Value = -1670.6667  # N22
# Row 31
$ws.Cells.Item(31, 8).Value = 2714.434  # H31
$ws.Cells.Item(31, 9).Value = 2063.5557  # I31
$ws.Cells.Item(31, 10).Value = 3390.3462  # J31
$ws.Cells.Item(31, 11).Value = 2063.5557  # K31
$ws.Cells.Item(31, 12).Value = 3390.3462  # L31
$ws.Cells.Item(31, 13).Value = -1768.5557  # M31
$ws.Cells.Item(31, 14).Value = -3980.3462  # N31
# Row 34
$ws.Cells.Item(34, 8).Value = 2714.434  # H34
$ws.Cells.Item(34, 9).Value = 2063.5557  # I34
$ws.Cells.Item(34, 10).Value = 3390.3462  # J34
$ws.Cells.Item(34, 11).Value = 2063.5557  # K34
$ws.Cells.Item(34, 12).Value = 3390.3462  # L34
$ws.Cells.Item(34, 13).Value = -1861.5557  # M34
$ws.Cells.Item(34, 14).Value = -3794.3462  # N34
# Row 58
$ws.Cells.Item(58, 8).Value = 1914.0303  # H58
$ws.Cells.Item(58, 9).Value = 2165  # I58
$ws.Cells.Item(58, 11).Value = 2165  # K58
$ws.Cells.Item(58, 13).Value = -1962  # M58
# Row 132
$ws.Cells.Item(132, 8).Value = 31880.53  # H132
$ws.Cells.Item(132, 9).Value = 47792.59  # I132
$ws.Cells.Item(132, 10).Value = 2708.4167  # J132
$ws.Cells.Item(132, 11).Value = 143377.77  # K132
$ws.Cells.Item(132, 12).Value = 8125.250100000001  # L132
$ws.Cells.Item(132, 13).Value = -140847.77  # M132
$ws.Cells.Item(132, 14).Value = -13185.2501  # N132
# Row 136
$ws.Cells.Item(136, 8).Value = 1914.0303  # H136
$ws.Cells.Item(136, 9).Value = 2165  # I136
$ws.Cells.Item(136, 11).Value = 6495  # K136
$ws.Cells.Item(136, 13).Value = -3945  # M136

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Cells.Item(4, 8).Value = 4923526.5  # H4
$ws.Cells.Item(4, 9).Value = 6444816.5  # I4
$ws.Cells.Item(4, 10).Value = 1500624.8  # J4
$ws.Cells.Item(4, 11).Value = 19334449.5  # K4
$ws.Cells.Item(4, 12).Value = 4501874.4  # L4
$ws.Cells.Item(4, 13).Value = -19334337.5  # M4
$ws.Cells.Item(4, 14).Value = -4502098.4  # N4
# Row 9
$ws.Cells.Item(9, 8).Value = 100000000  # H9
$ws.Cells.Item(9, 10).Value = 100000000  # J9
$ws.Cells.Item(9, 12).Value = 300000000  # L9
$ws.Cells.Item(9, 14).Value = -300000448  # N9
# Row 46
$ws.Cells.Item(46, 8).Value = 250005580  # H46
$ws.Cells.Item(46, 10).Value = 333338300  # J46
$ws.Cells.Item(46, 12).Value = 1000014900  # L46
$ws.Cells.Item(46, 14).Value = -1000015082  # N46
# Row 59
$ws.Cells.Item(59, 8).Value = 24000  # H59
$ws.Cells.Item(59, 9).Value = 0  # I59
$ws.Cells.Item(59, 10).Value = 24000  # J59
$ws.Cells.Item(59, 11).Value = 0  # K59
$ws.Cells.Item(59, 12).Value = 72000  # L59
$ws.Cells.Item(59, 13).ClearContents()  # M59
$ws.Cells.Item(59, 14).Value = -73080  # N59
# Row 64
$ws.Cells.Item(64, 8).Value = 9698.200000000001  # H64
$ws.Cells.Item(64, 10).Value = 9698.200000000001  # J64
$ws.Cells.Item(64, 12).Value = 29094.6  # L64
$ws.Cells.Item(64, 14).Value = -29634.6  # N64
# Row 67
$ws.Cells.Item(67, 8).Value = 9698.200000000001  # H67
$ws.Cells.Item(67, 10).Value = 9698.200000000001  # J67
$ws.Cells.Item(67, 12).Value = 29094.6  # L67
$ws.Cells.Item(67, 14).Value = -30966.6  # N67
# Row 81
$ws.Cells.Item(81, 8).Value = 30310974  # H81
$ws.Cells.Item(81, 9).Value = 111112590  # I81
$ws.Cells.Item(81, 10).Value = 10366.25  # J81
$ws.Cells.Item(81, 11).Value = 333337770  # K81
$ws.Cells.Item(81, 12).Value = 31098.75  # L81
$ws.Cells.Item(81, 13).Value = -333336647  # M81
$ws.Cells.Item(81, 14).Value = -33344.75  # N81
# Row 84
$ws.Cells.Item(84, 8).Value = 30310974  # H84
$ws.Cells.Item(84, 9).Value = 111112590  # I84
$ws.Cells.Item(84, 10).Value = 10366.25  # J84
$ws.Cells.Item(84, 11).Value = 1000013310  # K84
$ws.Cells.Item(84, 12).Value = 93296.25  # L84
$ws.Cells.Item(84, 13).Value = -1000007694  # M84
$ws.Cells.Item(84, 14).Value = -104528.25  # N84
# Row 92
$ws.Cells.Item(92, 8).Value = 566  # H92
$ws.Cells.Item(92, 9).Value = 566  # I92
$ws.Cells.Item(92, 10).Value = 0  # J92
$ws.Cells.Item(92, 11).Value = 1698  # K92
$ws.Cells.Item(92, 12).Value = 0  # L92
$ws.Cells.Item(92, 13).Value = -450  # M92
$ws.Cells.Item(92, 14).ClearContents()  # N92
# Row 115
$ws.Cells.Item(115, 8).Value = 1503.5  # H115
$ws.Cells.Item(115, 10).Value = 3699.25  # J115
$ws.Cells.Item(115, 12).Value = 11097.75  # L115
$ws.Cells.Item(115, 14).Value = -13447.75  # N115
# Row 121
$ws.Cells.Item(121, 8).Value = 3601.0833  # H121
$ws.Cells.Item(121, 10).Value = 9747.75  # J121
$ws.Cells.Item(121, 12).Value = 29243.25  # L121
$ws.Cells.Item(121, 14).Value = -31863.25  # N121
# Row 131
$ws.Cells.Item(131, 8).Value = 58626.168  # H131
$ws.Cells.Item(131, 10).Value = 103273  # J131
$ws.Cells.Item(131, 12).Value = 309819  # L131
$ws.Cells.Item(131, 14).Value = -319899  # N131
# Row 140
$ws.Cells.Item(140, 8).Value = 25643446  # H140
$ws.Cells.Item(140, 9).Value = 33335780  # I140
$ws.Cells.Item(140, 10).Value = 2333.3333  # J140
$ws.Cells.Item(140, 11).Value = 100007340  # K140
$ws.Cells.Item(140, 12).Value = 6999.999899999999  # L140
$ws.Cells.Item(140, 13).Value = -100002160  # M140
$ws.Cells.Item(140, 14).Value = -17359.9999  # N140

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80, 8).Value = 9738954  # H80
$ws.Cells.Item(80, 9).Value = 142979.34  # I80
$ws.Cells.Item(80, 10).Value = 45505770  # J80
$ws.Cells.Item(80, 11).Value = 142979.34  # K80
$ws.Cells.Item(80, 12).Value = 45505770  # L80
$ws.Cells.Item(80, 13).Value = -141981.34  # M80
$ws.Cells.Item(80, 14).Value = -45507766  # N80
# Row 83
$ws.Cells.Item(83, 8).Value = 9738954  # H83
$ws.Cells.Item(83, 9).Value = 142979.34  # I83
$ws.Cells.Item(83, 10).Value = 45505770  # J83
$ws.Cells.Item(83, 11).Value = 714896.7  # K83
$ws.Cells.Item(83, 12).Value = 227528850  # L83
$ws.Cells.Item(83, 13).Value = -709904.7  # M83
$ws.Cells.Item(83, 14).Value = -227538834  # N83
# Row 92
$ws.Cells.Item(92, 8).Value = 12123.8125  # H92
$ws.Cells.Item(92, 10).Value = 12123.8125  # J92
$ws.Cells.Item(92, 12).Value = 12123.8125  # L92
$ws.Cells.Item(92, 14).Value = -15867.8125  # N92

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Cells.Item(16, 8).Value = 786  # H16
$ws.Cells.Item(16, 10).Value = 1176.6666  # J16
$ws.Cells.Item(16, 12).Value = 1176.6666  # L16
$ws.Cells.Item(16, 14).Value = -1516.6666  # N16

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 75
$ws.Cells.Item(75, 8).Value = 27999.285  # H75
$ws.Cells.Item(75, 9).Value = 19199  # I75
$ws.Cells.Item(75, 10).Value = 50000  # J75
$ws.Cells.Item(75, 11).Value = 19199  # K75
$ws.Cells.Item(75, 12).Value = 50000  # L75
$ws.Cells.Item(75, 13).Value = -18263  # M75
$ws.Cells.Item(75, 14).Value = -51872  # N75
# Row 78
$ws.Cells.Item(78, 8).Value = 27999.285  # H78
$ws.Cells.Item(78, 9).Value = 19199  # I78
$ws.Cells.Item(78, 10).Value = 50000  # J78
$ws.Cells.Item(78, 11).Value = 57597  # K78
$ws.Cells.Item(78, 12).Value = 150000  # L78
$ws.Cells.Item(78, 13).Value = -52917  # M78
$ws.Cells.Item(78, 14).Value = -159360  # N78
# Row 81
$ws.Cells.Item(81, 8).Value = 163810.33  # H81
$ws.Cells.Item(81, 9).Value = 2572.5  # I81
$ws.Cells.Item(81, 11).Value = 5145  # K81
$ws.Cells.Item(81, 13).Value = -4084  # M81
# Row 84
$ws.Cells.Item(84, 8).Value = 163810.33  # H84
$ws.Cells.Item(84, 9).Value = 2572.5  # I84
$ws.Cells.Item(84, 11).Value = 25725  # K84
$ws.Cells.Item(84, 13).Value = -20421  # M84
# Row 132
$ws.Cells.Item(132, 8).Value = 1973.6167  # H132
$ws.Cells.Item(132, 9).Value = 1489.921  # I132
$ws.Cells.Item(132, 11).Value = 4469.763  # K132
$ws.Cells.Item(132, 13).Value = -1939.763  # M132
